$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1231.0869
$ws.Range("I17").Value = 770
$ws.Range("J17").Value = 1300.25
$ws.Range("K17").Value = 2310
$ws.Range("L17").Value = 3900.75
$ws.Range("M17").Value = -2142
$ws.Range("N17").Value = -4236.75
$ws.Range("H40").Value = 252624.75
$ws.Range("J40").Value = 502999.5
$ws.Range("L40").Value = 502999.5
$ws.Range("N40").Value = -503349.5
$ws.Range("H80").Value = 47136.453
$ws.Range("I80").Value = 50250
$ws.Range("J80").Value = 46444.555
$ws.Range("K80").Value = 150750
$ws.Range("L80").Value = 139333.665
$ws.Range("M80").Value = -149752
$ws.Range("N80").Value = -141329.665
$ws.Range("H83").Value = 47136.453
$ws.Range("I83").Value = 50250
$ws.Range("J83").Value = 46444.555
$ws.Range("K83").Value = 452250
$ws.Range("L83").Value = 418000.995
$ws.Range("M83").Value = -447258
$ws.Range("N83").Value = -427984.995
$ws.Range("H136").Value = 60520
$ws.Range("J136").Value = 60520
$ws.Range("L136").Value = 60520
$ws.Range("N136").Value = -70720

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 32931.758
$ws.Range("I74").Value = 47025
$ws.Range("K74").Value = 47025
$ws.Range("M74").Value = -46151
$ws.Range("H77").Value = 32931.758
$ws.Range("I77").Value = 47025
$ws.Range("K77").Value = 235125
$ws.Range("M77").Value = -230757
$ws.Range("H110").Value = 37038230
$ws.Range("I110").Value = 1304.2
$ws.Range("J110").Value = 83334390
$ws.Range("K110").Value = 1304.2
$ws.Range("L110").Value = 83334390
$ws.Range("M110").Value = 740.8
$ws.Range("N110").Value = -83338480
$ws.Range("H132").Value = 5252.185
$ws.Range("I132").Value = 4608.6113
$ws.Range("K132").Value = 13825.8339
$ws.Range("M132").Value = -11295.8339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5322202.5
$ws.Range("I134").Value = 7144252.5
$ws.Range("K134").Value = 21432757.5
$ws.Range("M134").Value = -21430222.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 43000
$ws.Range("J74").Value = 43000
$ws.Range("L74").Value = 43000
$ws.Range("N74").Value = -44748
$ws.Range("H77").Value = 43000
$ws.Range("J77").Value = 43000
$ws.Range("L77").Value = 129000
$ws.Range("N77").Value = -137736
$ws.Range("H99").Value = 4007.5264
$ws.Range("I99").Value = 2633
$ws.Range("K99").Value = 2633
$ws.Range("M99").Value = -1135
$ws.Range("H126").Value = 4007.5264
$ws.Range("I126").Value = 2633
$ws.Range("K126").Value = 7899
$ws.Range("M126").Value = -5429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1267.9375
$ws.Range("I5").Value = 833.1
$ws.Range("J5").Value = 1992.6666
$ws.Range("K5").Value = 2499.3
$ws.Range("L5").Value = 5977.9998
$ws.Range("M5").Value = -2387.3
$ws.Range("N5").Value = -6201.9998
$ws.Range("H12").Value = 478.52942
$ws.Range("J12").Value = 56.53846
$ws.Range("L12").Value = 169.61538
$ws.Range("N12").Value = -515.61538
$ws.Range("H68").Value = 33336288
$ws.Range("I68").Value = 33334106
$ws.Range("J68").Value = 33338472
$ws.Range("K68").Value = 100002318
$ws.Range("L68").Value = 100015416
$ws.Range("M68").Value = -100001507
$ws.Range("N68").Value = -100017038
$ws.Range("H71").Value = 33336288
$ws.Range("I71").Value = 33334106
$ws.Range("J71").Value = 33338472
$ws.Range("K71").Value = 300006954
$ws.Range("L71").Value = 300046248
$ws.Range("M71").Value = -300002898
$ws.Range("N71").Value = -300054360
$ws.Range("H97").Value = 416.5
$ws.Range("I97").Value = 411.14285
$ws.Range("J97").Value = 454
$ws.Range("K97").Value = 1233.42855
$ws.Range("L97").Value = 1362
$ws.Range("M97").Value = -737.4285500000001
$ws.Range("N97").Value = -2354
$ws.Range("H107").Value = 10000768
$ws.Range("I107").Value = 516.6667
$ws.Range("K107").Value = 1550.0001
$ws.Range("M107").Value = 369.9999
$ws.Range("H113").Value = 4626.3335
$ws.Range("I113").Value = 575
$ws.Range("J113").Value = 6099.5454
$ws.Range("K113").Value = 1725
$ws.Range("L113").Value = 18298.6362
$ws.Range("M113").Value = 445
$ws.Range("N113").Value = -22638.6362
$ws.Range("H114").Value = 47620348
$ws.Range("J114").Value = 66668290
$ws.Range("L114").Value = 200004870
$ws.Range("N114").Value = -200011378
$ws.Range("H132").Value = 5213.1943
$ws.Range("I132").Value = 2196.8
$ws.Range("J132").Value = 7367.7617
$ws.Range("K132").Value = 19771.2
$ws.Range("L132").Value = 66309.8553
$ws.Range("M132").Value = -17241.2
$ws.Range("N132").Value = -71369.8553
$ws.Range("H135").Value = 1267.9375
$ws.Range("I135").Value = 833.1
$ws.Range("J135").Value = 1992.6666
$ws.Range("K135").Value = 7497.900000000001
$ws.Range("L135").Value = 17933.9994
$ws.Range("M135").Value = -4962.900000000001
$ws.Range("N135").Value = -23003.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1770.9062
$ws.Range("I132").Value = 1522.4667
$ws.Range("K132").Value = 4567.4001
$ws.Range("M132").Value = -2037.4001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2707.8667
$ws.Range("I40").Value = 2739.8462
$ws.Range("K40").Value = 2739.8462
$ws.Range("M40").Value = -2603.8462
$ws.Range("H68").Value = 5681.125
$ws.Range("J68").Value = 5528.4287
$ws.Range("L68").Value = 5528.4287
$ws.Range("N68").Value = -7026.4287
$ws.Range("H71").Value = 5681.125
$ws.Range("J71").Value = 5528.4287
$ws.Range("L71").Value = 27642.1435
$ws.Range("N71").Value = -35130.14350000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14047321
$ws.Range("I81").Value = 946256.4399999999
$ws.Range("J81").Value = 50075250
$ws.Range("K81").Value = 1892512.88
$ws.Range("L81").Value = 100150500
$ws.Range("M81").Value = -1891451.88
$ws.Range("N81").Value = -100152622
$ws.Range("H84").Value = 14047321
$ws.Range("I84").Value = 946256.4399999999
$ws.Range("J84").Value = 50075250
$ws.Range("K84").Value = 9462564.399999999
$ws.Range("L84").Value = 500752500
$ws.Range("M84").Value = -9457260.399999999
$ws.Range("N84").Value = -500763108
$ws.Range("H122").Value = 2577.6775
$ws.Range("I122").Value = 1840.7693
$ws.Range("K122").Value = 5522.3079
$ws.Range("M122").Value = -3072.3079
$ws.Range("H136").Value = 23259586
$ws.Range("I136").Value = 45455740
$ws.Range("J136").Value = 6471.619
$ws.Range("K136").Value = 136367220
$ws.Range("L136").Value = 19414.857
$ws.Range("M136").Value = -136364670
$ws.Range("N136").Value = -24514.857
